# Refresh the "cryptos" price/volume table (rows 2-51) with the latest
# scraped values, as produced by the GitHub Actions bot on
# Thu Apr 6 09:20:44 UTC 2023. Row 13/14 also swap places (Solana now
# ranks above WrappedEther).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as literal text (not auto-converted numbers)
# so values like "1.000" or "27.996.44" keep their exact textual form.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.996.44'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '1.882.86'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '312.70'
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.4992'
$ws.Range('E7').Value = '  -3.62%  '
$ws.Range('D8').Value = '0.3856'
$ws.Range('E8').Value = '  -2.66%  '
$ws.Range('D9').Value = '0.09173'
$ws.Range('E9').Value = '  -5.68%  '
$ws.Range('D10').Value = '1.121'
$ws.Range('E10').Value = '  -2.88%  '
$ws.Range('D11').Value = '41.67'
$ws.Range('D12').Value = '6.331'
$ws.Range('E12').Value = '  -3.14%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '20.73'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.876.13'
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').Value = '7.284'
$ws.Range('E15').Value = '  -3.14%  '
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '0.00001105'
$ws.Range('E17').Value = '  -2.86%  '
$ws.Range('D18').Value = '91.34'
$ws.Range('E18').Value = '  -3.68%  '
$ws.Range('D19').Value = '0.06632'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('D20').Value = '17.98'
$ws.Range('E20').Value = '  -1.34%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = '6.176'
$ws.Range('E22').Value = '  -2.36%  '
$ws.Range('D23').Value = '28.044.42'
$ws.Range('E23').Value = '  -2.25%  '
$ws.Range('D24').Value = '11.39'
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('D25').Value = '2.302'
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('D26').Value = '2.100.61'
$ws.Range('E26').Value = '  -1.34%  '
$ws.Range('D27').Value = '2.543'
$ws.Range('E27').Value = '  -5.25%  '
$ws.Range('D28').Value = '157.58'
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('E29').Value = '  -2.55%  '
$ws.Range('D30').Value = '126.61'
$ws.Range('E30').Value = '  -1.85%  '
$ws.Range('D31').Value = '0.1056'
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('D32').Value = '1.066'
$ws.Range('E32').Value = '  -4.50%  '
$ws.Range('D33').Value = '5.588'
$ws.Range('E33').Value = '  -3.17%  '
$ws.Range('D34').Value = '3.580'
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('D35').Value = '9.334'
$ws.Range('E35').Value = '  -5.74%  '
$ws.Range('D36').Value = '0.06583'
$ws.Range('E36').Value = '  -3.16%  '
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('D38').Value = '0.2189'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('D39').Value = '1.287'
$ws.Range('E39').Value = '  +8.02%  '
$ws.Range('D40').Value = '1.207'
$ws.Range('E40').Value = '  -5.42%  '
$ws.Range('D41').Value = '0.6401'
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('D42').Value = '11.51'
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range('D43').Value = '4.930'
$ws.Range('E43').Value = '  -3.30%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = '13.30'
$ws.Range('E45').Value = '  -2.62%  '
$ws.Range('D46').Value = '0.6040'
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('D47').Value = '1.293'
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('D48').Value = '3.667'
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('D49').Value = '1.988'
$ws.Range('E49').Value = '  -2.50%  '
$ws.Range('D50').Value = '1.210'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').Value = '121.30'
$ws.Range('E51').Value = '  -3.05%  '

# Restore default (Normal) style so no stray number-format styling is introduced
$textRange.Style = "Normal"
